$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price and volume(1h) values
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("E3").Value = "  -2.96%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("E7").Value = "  -6.75%  "
$ws.Range("E8").Value = "  +14.88%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "3.096.65"
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("E11").Value = "  +9.09%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("E13").Value = "  -6.47%  "
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "88.854.64"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("E16").Value = "  -2.38%  "
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("D18").Value = "3.102.04"
$ws.Range("E18").Value = "  -3.55%  "
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("E20").Value = "  -6.51%  "
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("E23").Value = "  -3.26%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("E25").Value = "  +5.21%  "
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("D28").Value = "3.235.71"
$ws.Range("E28").Value = "  -4.15%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +8.20%  "
$ws.Range("E31").Value = "  +7.50%  "
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("E34").Value = "  -10.43%  "
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("E36").Value = "  -3.04%  "
$ws.Range("E37").Value = "  -4.99%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  +3.80%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  -2.45%  "
$ws.Range("E44").Value = "  -5.98%  "
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("E46").Value = "  +4.77%  "
$ws.Range("E47").Value = "  +13.86%  "
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("E49").Value = "  -6.59%  "
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("E51").Value = "  -5.26%  "
